{"js": "// Remove the `help(\"round\")` SourceCode paragraph that followed the\n// `## [1] 9.333333` output paragraph (author replaced the list-comprehension\n// question's old scratch snippet while experimenting with vector ops).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf('help(\"round\")') !== -1) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the `help(\"round\")` SourceCode paragraph that followed the\n# `## [1] 9.333333` output paragraph (author replaced the list-comprehension\n# question's old scratch snippet while experimenting with vector ops).\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text.TrimEnd(\"`r\")\n  if ($t -eq 'help(\"round\")') {\n    $p.Range.Delete()\n  }\n}\n"}
